# UPDATE technology portfolios for Norway
$wb = $excel.ActiveWorkbook

$ws2025 = $wb.Worksheets.Item("2025")
$ws2030 = $wb.Worksheets.Item("2030")
$ws2035 = $wb.Worksheets.Item("2035")
$ws2040 = $wb.Worksheets.Item("2040")
$ws2045 = $wb.Worksheets.Item("2045")
$ws2050 = $wb.Worksheets.Item("2050")

# Update the base investment cost for 2025; the dependent years recalc
# automatically since they reference '2025'!A2.
$ws2025.Range("A2").Value = 295000

# Update p_ieh_elec (column B) for every year from 11 to 10.
$ws2025.Range("B2").Value = 10
$ws2030.Range("B2").Value = 10
$ws2035.Range("B2").Value = 10
$ws2040.Range("B2").Value = 10
$ws2045.Range("B2").Value = 10
$ws2050.Range("B2").Value = 10

# Move the active tab selection from 2025 to 2050.
$ws2050.Activate()
